$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in B1 from "AccountNumber" to "Account Number"
$ws.Range("B1").Value = "Account Number"

# Update column B's width to reflect the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 15.5703125

# Move the active selection to B1 as recorded in the saved view state
$ws.Range("B1").Select() | Out-Null
